$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "58.918.69"
Set-TextValue 2 5 "  -3.25%  "

Set-TextValue 3 5 "  -3.85%  "

Set-TextValue 4 5 "  +0.00%  "

Set-TextValue 5 4 "540.52"
Set-TextValue 5 5 "  -4.43%  "

Set-TextValue 6 4 "136.27"
Set-TextValue 6 5 "  -8.20%  "

Set-TextValue 7 5 "  -0.05%  "

Set-TextValue 8 4 "3.227.78"
Set-TextValue 8 5 "  -3.96%  "

Set-TextValue 9 5 "  -4.09%  "

Set-TextValue 10 5 "  -4.49%  "

Set-TextValue 11 5 "  -5.53%  "

Set-TextValue 12 5 "  -4.40%  "

Set-TextValue 13 4 "3.782.96"
Set-TextValue 13 5 "  -3.90%  "

Set-TextValue 14 5 "  -1.00%  "

Set-TextValue 15 5 "  -7.12%  "

Set-TextValue 16 4 "3.227.84"
Set-TextValue 16 5 "  -4.32%  "

Set-TextValue 17 5 "  -5.82%  "

Set-TextValue 18 4 "58.951.17"
Set-TextValue 18 5 "  -3.36%  "

Set-TextValue 19 4 "5.90"
Set-TextValue 19 5 "  -6.93%  "

Set-TextValue 20 4 "13.39"
Set-TextValue 20 5 "  -5.62%  "

Set-TextValue 21 4 "8.24"
Set-TextValue 21 5 "  -6.73%  "

Set-TextValue 22 4 "362.82"
Set-TextValue 22 5 "  -3.00%  "

Set-TextValue 23 5 "  -0.12%  "

Set-TextValue 24 4 "70.56"
Set-TextValue 24 5 "  -6.28%  "

Set-TextValue 25 4 "0.521"
Set-TextValue 25 5 "  -6.63%  "

Set-TextValue 26 4 "3.363.15"
Set-TextValue 26 5 "  -3.94%  "

Set-TextValue 27 5 "  -2.64%  "

Set-TextValue 28 4 "0.0₃0972"
Set-TextValue 28 5 "  -10.01%  "

Set-TextValue 29 4 "0.997"
Set-TextValue 29 5 "  -0.57%  "

Set-TextValue 30 4 "7.09"
Set-TextValue 30 5 "  -3.84%  "

Set-TextValue 31 5 "  -0.07%  "

Set-TextValue 32 4 "1.94"
Set-TextValue 32 5 "  -6.52%  "

Set-TextValue 33 5 "  -7.51%  "

Set-TextValue 34 4 "21.96"
Set-TextValue 34 5 "  -3.80%  "

Set-TextValue 36 5 "  -7.64%  "

Set-TextValue 37 4 "162.69"
Set-TextValue 37 5 "  -4.17%  "

Set-TextValue 38 4 "6.43"
Set-TextValue 38 5 "  -5.10%  "

Set-TextValue 39 5 "  -6.03%  "

Set-TextValue 40 4 "26.29"
Set-TextValue 40 5 "  -9.63%  "

Set-TextValue 41 5 "  -4.65%  "

Set-TextValue 42 4 "3.258.96"
Set-TextValue 42 5 "  -4.04%  "

Set-TextValue 43 4 "41.07"
Set-TextValue 43 5 "  -2.92%  "

Set-TextValue 44 5 "  -5.65%  "

Set-TextValue 45 5 "  -3.70%  "

Set-TextValue 46 4 "4.04"
Set-TextValue 46 5 "  -5.57%  "

Set-TextValue 47 5 "  -6.13%  "

Set-TextValue 48 5 "  -0.11%  "

Set-TextValue 49 4 "2.303.28"
Set-TextValue 49 5 "  -7.39%  "

Set-TextValue 50 5 "  -5.03%  "

Set-TextValue 51 4 "20.91"
Set-TextValue 51 5 "  -6.69%  "
